$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("Bradley Beal", "PG,SG,SF", "Phoenix Suns"),
    @("Malik Beasley", "SG", "Detroit Pistons"),
    @("Dyson Daniels", "PG,SG", "Atlanta Hawks"),
    @("De'Andre Hunter", "SF,PF", "Atlanta Hawks"),
    @("Alexandre Sarr", "PF,C", "Washington Wizards"),
    @("Donovan Mitchell", "PG,SG", "Cleveland Cavaliers"),
    @("Payton Pritchard", "PG", "Boston Celtics"),
    @("Josh Hart", "SF,PF", "New York Knicks"),
    @("Myles Turner", "C", "Indiana Pacers"),
    @("Domantas Sabonis", "C", "Sacramento Kings"),
    @("Deandre Ayton", "C", "Portland Trail Blazers"),
    @("Tari Eason", "SF,PF", "Houston Rockets"),
    @("Victor Wembanyama", "C", "San Antonio Spurs"),
    @("Jamal Murray", "PG,SG", "Denver Nuggets"),
    @("Michael Porter Jr.", "SF,PF", "Denver Nuggets"),
    @("Zion Williamson", "PF,C", "New Orleans Pelicans"),
    @("Kristaps Porzingis", "PF,C", "Boston Celtics")
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
    $ws.Cells.Item($row, 3).Value = $data[$i][2]
}
